# Plan van Aanpak - insert a new paragraph "Testen.?" right after the
# "-Comments (Roos en Lubbert)" paragraph (before the page-break run that
# precedes "Benodigdheden"), matching the run formatting (szCs=17) used
# by the rest of that paragraph.

$d = $word.ActiveDocument

# Locate the unique "(Roos en Lubbert)" run that ends the "-Comments" line.
$srcRange = $d.Content
$searchText = "Roos en Lubbert)"
$found = $srcRange.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text '$searchText'"
}

# Grab a copy of this run's full formatting (FormattedText acts like a
# format-painter snapshot: text + run properties) so the new paragraph's
# run ends up with the exact same <w:rPr> (szCs=17, nothing else).
$srcFormatted = $srcRange.FormattedText
$srcLen = $srcRange.End - $srcRange.Start

# Split the paragraph right after "(Roos en Lubbert)": this closes the
# current paragraph and opens a new (empty) one that inherits the
# paragraph mark's run properties, in front of the existing page-break
# run ("Benodigdheden" page).
$srcRange.Collapse(0)
$srcRange.InsertParagraphAfter()

# Position just inside the freshly created paragraph (right after the
# paragraph-mark character that was just inserted).
$newParaStart = $srcRange.End + 1

# Paste the captured formatting (creates its own run carrying the same
# <w:rPr> as the source, decoupled from the following Arial run).
$pasteRange = $d.Range($newParaStart, $newParaStart)
$pasteRange.FormattedText = $srcFormatted

# The pasted run currently holds a copy of "Roos en Lubbert)" text;
# swap its text for "Testen.?" while keeping the run's formatting.
$newRunRange = $d.Range($newParaStart, $newParaStart + $srcLen)
$newRunRange.Text = "Testen.?"

Write-Output "Inserted 'Testen.?' paragraph."
